# Generated edit script for localization-status.xlsx
# Commit message: Generate Report for Handoff
$wb = $excel.ActiveWorkbook

# ===== Sheet: Overview =====
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A1").Hyperlinks.Delete()  # clear all existing hyperlinks on this sheet

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/calleeMd1.md", "", "", "calleeMd1.md")
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-23 11:10:00"

$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/calleeMd2.md", "", "", "calleeMd2.md")
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-23 11:10:00"

$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/callerMd1.md", "", "", "callerMd1.md")
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-23 11:10:00"

$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/callerMd2.md", "", "", "callerMd2.md")
$ws1.Range("B5").Value = "Ready for handoff"
$ws1.Range("C5").Value = "Ready for handoff"
$ws1.Range("D5").Value = "2016-03-23 11:10:00"

# ===== Sheet: zh-cn =====
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()  # clear all existing hyperlinks on this sheet

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/calleeMd1.md", "", "", "calleeMd1.md")
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9aba06ed984db768f4fd99a840afcba83efdc2a8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.zh-cn.xlf")
$ws2.Range("E2").Value = "2016-03-23 11:09:56"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("J2").Value = "Include"
$ws2.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/calleeMd2.md", "", "", "calleeMd2.md")
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9aba06ed984db768f4fd99a840afcba83efdc2a8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.zh-cn.xlf")
$ws2.Range("E3").Value = "2016-03-23 11:09:56"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("J3").Value = "Include"
$ws2.Range("K3").Value = "e2e\callerMd1.md"

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/callerMd1.md", "", "", "callerMd1.md")
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9aba06ed984db768f4fd99a840afcba83efdc2a8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.zh-cn.xlf")
$ws2.Range("E4").Value = "2016-03-23 11:09:56"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("J4").Value = "Include"
$ws2.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/callerMd2.md", "", "", "callerMd2.md")
$ws2.Range("B5").Value = ".md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9aba06ed984db768f4fd99a840afcba83efdc2a8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.zh-cn.xlf")
$ws2.Range("E5").Value = "2016-03-23 11:09:56"
$ws2.Range("H5").Value = "0001-01-01 00:00:00"
$ws2.Range("J5").Value = "Include"
$ws2.Range("I5").Value = "e2e\calleeMd1.md"

# ===== Sheet: de-de =====
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()  # clear all existing hyperlinks on this sheet

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/calleeMd1.md", "", "", "calleeMd1.md")
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5776397db0ccb90c7738b160b9fe010e431a29c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf", "", "", "calleeMd1.e8f5ecec2b522eb147a4ff0ca19ca72e17f2186d.de-de.xlf")
$ws3.Range("E2").Value = "2016-03-23 11:10:00"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("J2").Value = "Include"
$ws3.Range("K2").Value = "e2e\callerMd2.md,`ne2e\callerMd1.md"

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/calleeMd2.md", "", "", "calleeMd2.md")
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5776397db0ccb90c7738b160b9fe010e431a29c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf", "", "", "calleeMd2.63b76063f058ecc63ff1dda71ea2a67db72ae6e1.de-de.xlf")
$ws3.Range("E3").Value = "2016-03-23 11:10:00"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("J3").Value = "Include"
$ws3.Range("K3").Value = "e2e\callerMd1.md"

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/callerMd1.md", "", "", "callerMd1.md")
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5776397db0ccb90c7738b160b9fe010e431a29c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf", "", "", "callerMd1.a3bf9f4e7fa2750ec06df0b78a76ae5cafa0e0fd.de-de.xlf")
$ws3.Range("E4").Value = "2016-03-23 11:10:00"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("J4").Value = "Include"
$ws3.Range("I4").Value = "e2e\calleeMd1.md,`ne2e\calleeMd2.md"

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d278f271245a135d8c04667de7234ee3d2d93d51/e2e/callerMd2.md", "", "", "callerMd2.md")
$ws3.Range("B5").Value = ".md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5776397db0ccb90c7738b160b9fe010e431a29c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf", "", "", "callerMd2.c7d976edeb9cd5406eae7aba4c05d6d92e81ae95.de-de.xlf")
$ws3.Range("E5").Value = "2016-03-23 11:10:00"
$ws3.Range("H5").Value = "0001-01-01 00:00:00"
$ws3.Range("J5").Value = "Include"
$ws3.Range("I5").Value = "e2e\calleeMd1.md"
